# Update the "indirect-expenses" financial analysis sheet with refreshed figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# As-of date (B1) moved from 01-May-2024 to 01-Jan-2025
$ws.Range("B1").Value = 45658

# Updated expense figures
$ws.Range("B2").Value = 532489.14
$ws.Range("B3").Value = 72335.17
$ws.Range("B7").Value = 310
$ws.Range("B9").Value = 30000
$ws.Range("B10").Value = 3468
$ws.Range("B13").Value = 9554.82

# B14 (General Repairs & Maintanance) now has a value
$ws.Range("B14").Value = 5100

$ws.Range("B20").Value = 4920
$ws.Range("B21").Value = 3980.85
$ws.Range("B22").Value = 28260
$ws.Range("B23").Value = 80500

# B27 (Repairs & Service Charges) cleared out
$ws.Range("B27").ClearContents()

$ws.Range("B29").Value = 8.5500000000000007

# B30 (Seminar, Training, Dvpt. Exp.) cleared out
$ws.Range("B30").ClearContents()

$ws.Range("B31").Value = 25566.97

# B32 (Service charges) now has a value
$ws.Range("B32").Value = 24483

$ws.Range("B33").Value = 4609
$ws.Range("B35").Value = 6731.78
$ws.Range("B36").Value = 5351
$ws.Range("B38").Value = 2875
$ws.Range("B39").Value = 192935
$ws.Range("B40").Value = 944632.41
$ws.Range("B41").Value = 10643.6
$ws.Range("B42").Value = 193459.81
$ws.Range("B43").Value = 97333
$ws.Range("B45").Value = 524196

# B47 (Interest On Term Loan) cleared out
$ws.Range("B47").ClearContents()

$ws.Range("B51").Value = 53840.53

# B54 (Gift Articles) cleared out
$ws.Range("B54").ClearContents()

$ws.Range("B55").Value = 30.53
$ws.Range("B56").Value = 41600

# B58 (Discount Allowed) cleared out
$ws.Range("B58").ClearContents()

$ws.Range("B60").Value = 12210
$ws.Range("B64").Value = 1530962.08
$ws.Range("B66").Value = 532489
$ws.Range("B67").Value = 53841
$ws.Range("B68").Value = 586330
$ws.Range("B69").Value = 728299
$ws.Range("B71").Value = 0
$ws.Range("B72").Value = 216333
$ws.Range("B73").Value = 944632
$ws.Range("B74").Value = 1530962

# Update the active selection to K15 (as recorded in the saved view state)
[void]$ws.Range("K15").Select()
